$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Area" / "Subarea1" / "Subarea2" sample values so they contain
# special (accented) characters, exercising the authors-import bug fix
# described in the commit message.
$ws.Range("M2").Value = "congreso de la sociedad venezolana de física"
$ws.Range("N2").Value = "física nuclear"
$ws.Range("M3").Value = "ciencias sociales"
$ws.Range("N3").Value = "educación"
$ws.Range("O3").Value = "ciencias sociales"

# Extend the sheet's used range out to column S (an extra blank column was
# left behind in the saved file). Write then clear a value so the cell is
# materialised, then restore its default (unstyled) formatting by pasting
# the formatting of an already-default-styled cell onto it.
$ws.Range("S1").Value = "x"
$ws.Range("S1").ClearContents()
$ws.Range("A1").Copy()
[void]$ws.Range("S1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Match the final cursor/selection state recorded in the saved file.
[void]$ws.Range("R1:R3").Select()

Write-Output "done"
